$wb = $excel.ActiveWorkbook

# --- Create the new "US8" sheet by duplicating "US7" (same layout/styles), ---
# --- placed immediately before "Issues" (i.e. right after "US7").          ---
$sheetUS7 = $wb.Worksheets.Item("US7")
$sheetIssues = $wb.Worksheets.Item("Issues")
$sheetUS7.Copy($sheetIssues)
$newSheet = $wb.Worksheets.Item("US7 (2)")
$newSheet.Name = "US8"

# Remove the old US7 task rows (9-19) that don't apply to US8.
$newSheet.Range("B9:F19").EntireRow.Delete()

# --- Populate new content. Cells containing brand-new text are written in ---
# --- this specific order so shared-string indices come out 104..107 in   ---
# --- the same order as the authoritative edit.                           ---
$newSheet.Range("C6").Value = "Calculate Total working hours and total deficit/overtime hours for specified date range"
$newSheet.Range("C8").Value = "write test cases"
$newSheet.Range("B3").Value = "US 8"
$newSheet.Range("C7").Value = "UI - show total working hours and deficit/overtime hours "

# Remaining (already-existing) strings / values.
$newSheet.Range("C3").Value = "Show Total hrs for the selected data"

$newSheet.Range("B6").Value = 1
$newSheet.Range("D6").Value = 2
$newSheet.Range("E6").Value = "Sanket"
$newSheet.Range("F6").Value = "To do"

$newSheet.Range("B7").Value = 2
$newSheet.Range("D7").Value = 2
$newSheet.Range("E7").Value = "Shweta"
$newSheet.Range("F7").Value = "To do"

$newSheet.Range("B8").Value = 3
$newSheet.Range("D8").Value = 2
$newSheet.Range("E8").Value = "Sidhdesh"
$newSheet.Range("F8").Value = "To do"

# Widen column C to fit the longer task descriptions, matching the edited file.
$newSheet.Columns.Item(3).ColumnWidth = 79.140625
$newSheet.Range("C7").Select()

# --- Rename old "US7" sheet to "US7 " (trailing space), as in the edit. ---
$sheetUS7.Name = "US7 "

# --- Make the new sheet the active / selected tab. ---
$newSheet.Activate()
